$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary_counts")

# Insert a new row above row 13, shifting rows 13-14 down to 14-15.
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "Number of events with both any university response coding and any police coding"
$ws.Range("B13").Value = 360
